$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1203.5
$ws.Range("I6").Value = 919.1539
$ws.Range("K6").Value = 2757.4617
$ws.Range("M6").Value = -2645.4617
# Row 17
$ws.Range("H17").Value = 1621.4
$ws.Range("J17").Value = 1690.4445
$ws.Range("L17").Value = 5071.333500000001
$ws.Range("N17").Value = -5407.333500000001
# Row 33
$ws.Range("H33").Value = 188.62962
$ws.Range("I33").Value = 105.63158
$ws.Range("K33").Value = 105.63158
$ws.Range("M33").Value = 123.36842
# Row 39
$ws.Range("H39").Value = 495.33334
$ws.Range("I39").Value = 495.33334
$ws.Range("K39").Value = 1486.00002
$ws.Range("M39").Value = -1190.00002
# Row 40
$ws.Range("H40").Value = 3160
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5350
# Row 64
$ws.Range("H64").Value = 25560.1
$ws.Range("I64").Value = 34028.715
$ws.Range("J64").Value = 5800
$ws.Range("K64").Value = 34028.715
$ws.Range("L64").Value = 5800
$ws.Range("M64").Value = -33780.715
$ws.Range("N64").Value = -6296
# Row 67
$ws.Range("H67").Value = 25560.1
$ws.Range("I67").Value = 34028.715
$ws.Range("J67").Value = 5800
$ws.Range("K67").Value = 34028.715
$ws.Range("L67").Value = 5800
$ws.Range("M67").Value = -33170.715
$ws.Range("N67").Value = -7516
# Row 111
$ws.Range("H111").Value = 8306.4
$ws.Range("J111").Value = 10516
$ws.Range("L111").Value = 31548
$ws.Range("N111").Value = -37682
# Row 115
$ws.Range("H115").Value = 933.3333
$ws.Range("I115").Value = 933.3333
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 2799.9999
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -1232.9999
$ws.Range("N115").ClearContents()
# Row 129
$ws.Range("H129").Value = 1992.2727
$ws.Range("I129").Value = 1490.125
$ws.Range("K129").Value = 4470.375
$ws.Range("M129").Value = 529.625
# Row 132
$ws.Range("H132").Value = 22245
$ws.Range("I132").Value = 25621.725
$ws.Range("K132").Value = 76865.17499999999
$ws.Range("M132").Value = -74335.17499999999
# Row 137
$ws.Range("H137").Value = 25942.334
$ws.Range("I137").Value = 75749.25
$ws.Range("J137").Value = 1038.875
$ws.Range("K137").Value = 227247.75
$ws.Range("L137").Value = 3116.625
$ws.Range("M137").Value = -224697.75
$ws.Range("N137").Value = -8216.625
# Row 138
$ws.Range("H138").Value = 35872.633
$ws.Range("I138").Value = 2227.158
$ws.Range("K138").Value = 6681.474
$ws.Range("M138").Value = -1541.474
# Row 141
$ws.Range("H141").Value = 3090.4285
$ws.Range("I141").Value = 3090.4285
$ws.Range("K141").Value = 9271.2855
$ws.Range("M141").Value = -4091.2855

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3796.5908
$ws.Range("I2").Value = 3727.0527
$ws.Range("K2").Value = 3727.0527
$ws.Range("M2").Value = -3614.0527
# Row 25
$ws.Range("H25").Value = 20600
$ws.Range("J25").Value = 21700
$ws.Range("L25").Value = 21700
$ws.Range("N25").Value = -22504
# Row 45
$ws.Range("H45").Value = 4124.5713
$ws.Range("I45").Value = 2250.6667
$ws.Range("K45").Value = 2250.6667
$ws.Range("M45").Value = -1873.6667
# Row 61
$ws.Range("H61").Value = 5120.2964
$ws.Range("I61").Value = 1101.5454
$ws.Range("K61").Value = 1101.5454
$ws.Range("M61").Value = -889.5454
# Row 74
$ws.Range("H74").Value = 437576.56
$ws.Range("I74").Value = 750826.4
$ws.Range("K74").Value = 750826.4
$ws.Range("M74").Value = -749952.4
# Row 77
$ws.Range("H77").Value = 437576.56
$ws.Range("I77").Value = 750826.4
$ws.Range("K77").Value = 3754132
$ws.Range("M77").Value = -3749764
# Row 95
$ws.Range("H95").Value = 58999
$ws.Range("I95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("M95").ClearContents()
# Row 110
$ws.Range("H110").Value = 6927.4194
$ws.Range("I110").Value = 6927.4194
$ws.Range("K110").Value = 6927.4194
$ws.Range("M110").Value = -4882.4194
# Row 116
$ws.Range("H116").Value = 3796.5908
$ws.Range("I116").Value = 3727.0527
$ws.Range("K116").Value = 3727.0527
$ws.Range("M116").Value = -1433.0527
# Row 136
$ws.Range("H136").Value = 5120.2964
$ws.Range("I136").Value = 1101.5454
$ws.Range("K136").Value = 3304.6362
$ws.Range("M136").Value = -754.6361999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3796.5908
$ws.Range("I3").Value = 3727.0527
$ws.Range("K3").Value = 3727.0527
$ws.Range("M3").Value = -3613.0527
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1633.3334
$ws.Range("I16").Value = 1633.3334
$ws.Range("K16").Value = 1633.3334
$ws.Range("M16").Value = -1346.3334
# Row 31
$ws.Range("H31").Value = 2780336.2
$ws.Range("I31").Value = 5556451.5
$ws.Range("J31").Value = 4220.722
$ws.Range("K31").Value = 5556451.5
$ws.Range("L31").Value = 4220.722
$ws.Range("M31").Value = -5556156.5
$ws.Range("N31").Value = -4810.722
# Row 34
$ws.Range("H34").Value = 2780336.2
$ws.Range("I34").Value = 5556451.5
$ws.Range("J34").Value = 4220.722
$ws.Range("K34").Value = 5556451.5
$ws.Range("L34").Value = 4220.722
$ws.Range("M34").Value = -5556249.5
$ws.Range("N34").Value = -4624.722
# Row 113
$ws.Range("H113").Value = 1633.3334
$ws.Range("I113").Value = 1633.3334
$ws.Range("K113").Value = 1633.3334
$ws.Range("M113").Value = 536.6666
# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
# Row 132
$ws.Range("H132").Value = 127980.5
$ws.Range("I132").Value = 334003.66
$ws.Range("K132").Value = 1002010.98
$ws.Range("M132").Value = -999480.98

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 100925864
$ws.Range("I4").Value = 126403870
$ws.Range("K4").Value = 379211610
$ws.Range("M4").Value = -379211498
# Row 70
$ws.Range("H70").Value = 4258.5713
# Row 73
$ws.Range("H73").Value = 4258.5713
# Row 115
$ws.Range("H115").Value = 299
$ws.Range("J115").Value = 299
$ws.Range("L115").Value = 897
$ws.Range("N115").Value = -3247

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2395.56
$ws.Range("J102").Value = 2708.5454
$ws.Range("L102").Value = 2708.5454
$ws.Range("N102").Value = -5952.5454

$ws = $wb.Worksheets.Item("LTW")
# Row 43
$ws.Range("H43").Value = 18752.5
$ws.Range("J43").Value = 18752.5
$ws.Range("L43").Value = 18752.5
$ws.Range("N43").Value = -19138.5
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
# Row 122
$ws.Range("H122").Value = 3078.7334
$ws.Range("I122").Value = 2999.125
$ws.Range("K122").Value = 8997.375
$ws.Range("M122").Value = -6547.375
# Row 132
$ws.Range("H132").Value = 4997.5
$ws.Range("I132").Value = 4997.5
$ws.Range("K132").Value = 14992.5
$ws.Range("M132").Value = -12462.5

$ws = $wb.Worksheets.Item("WVR")
# Row 43
$ws.Range("H43").Value = 37999
$ws.Range("I43").Value = 37999
$ws.Range("K43").Value = 37999
$ws.Range("M43").Value = -37850
# Row 107
$ws.Range("H107").Value = 1247
$ws.Range("J107").Value = 943.6667
$ws.Range("L107").Value = 2831.0001
$ws.Range("N107").Value = -6671.0001
# Row 122
$ws.Range("H122").Value = 54049.31
$ws.Range("I122").Value = 60643.434
$ws.Range("J122").Value = 3494.3333
$ws.Range("K122").Value = 181930.302
$ws.Range("L122").Value = 10482.9999
$ws.Range("M122").Value = -179480.302
$ws.Range("N122").Value = -15382.9999
# Row 136
$ws.Range("H136").Value = 22750.484
$ws.Range("I136").Value = 25306.111
$ws.Range("J136").Value = 5500
$ws.Range("K136").Value = 75918.333
$ws.Range("L136").Value = 16500
$ws.Range("M136").Value = -73368.333
$ws.Range("N136").Value = -21600
